# Issue #31 Repackage pages and pageComponents should be separate
#
# Adds six new rows (26-31) to the "Issues" sheet log, and bumps the
# Priority on row 25 from 1 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Row 25 priority changed from 1 to 2
$ws.Cells.Item(25, 2).Value = 2

# Row 26: error / Better error handling when REST server not available
$ws.Cells.Item(26, 1).Value = 26
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(26, 4).Value = "error"
$ws.Cells.Item(26, 5).Value = "Better error handling when REST server not available"

# Row 27: Hearbeat / Requires 24 / Settings option ...
$ws.Cells.Item(27, 1).Value = 27
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(27, 5).Value = "Hearbeat"
$ws.Cells.Item(27, 6).Value = 24
$ws.Cells.Item(27, 8).Value = "Settings option to allow for connectivity to REST server to be established"

# Row 28: icons same size / Should make the tick/cross icons the same size ...
$ws.Cells.Item(28, 1).Value = 28
$ws.Cells.Item(28, 2).Value = 2
$ws.Cells.Item(28, 5).Value = "icons same size"
$ws.Cells.Item(28, 8).Value = "Should make the tick/cross icons the same size so they do not move the screen when they are changeg"

# Row 29: Should add return from browse playlist to main menu
$ws.Cells.Item(29, 1).Value = 29
$ws.Cells.Item(29, 2).Value = 2
$ws.Cells.Item(29, 5).Value = "Should add return from browse playlist to main menu"

# Row 30: Should have an return nav bar on add playlist
$ws.Cells.Item(30, 1).Value = 30
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(30, 5).Value = "Should have an return nav bar on add playlist"

# Row 31: arch / Repackage pages and pageComponents should be separate
$ws.Cells.Item(31, 1).Value = 31
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = "DONE"
$ws.Cells.Item(31, 4).Value = "arch"
$ws.Cells.Item(31, 5).Value = "Repackage pages and pageComponents should be separate"

# Row heights Excel auto-calculated for the wrapped text in the new rows
$ws.Rows.Item(26).RowHeight = 43.5
$ws.Rows.Item(27).RowHeight = 29
$ws.Rows.Item(28).RowHeight = 29
$ws.Rows.Item(29).RowHeight = 43.5
$ws.Rows.Item(30).RowHeight = 29
$ws.Rows.Item(31).RowHeight = 43.5

# Selection follows the last edited cell
$ws.Range("E31").Select() | Out-Null
